# Atualização dos notebooks para inclusão da inf de Renda Mensal
#
# Adds a new worksheet "DIME_RENDA_MENSAL" at the end of the workbook,
# populated with the QUESTIONARIO_ID / QUESTIONARIO_DESCRICAO lookup table
# for monthly income brackets, mirroring the layout of the other DIME_*
# lookup sheets already in the workbook.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (DIME_ESTADO_CIVIL)
# so it becomes the final tab, then rename + activate it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "DIME_RENDA_MENSAL"

# Header row
$ws.Range("A1").Value = "QUESTIONARIO_ID"
$ws.Range("B1").Value = "QUESTIONARIO_DESCRICAO"

# Income-bracket lookup rows (questionnaire options A..Q), entered the same
# way the author built the sheet: column A for the first block of rows,
# then column B for that same block, then the remaining rows filled in one
# at a time (A then B per row).
$ws.Range("A2").Value = "A"
$ws.Range("A3").Value = "B"
$ws.Range("A4").Value = "C"
$ws.Range("A5").Value = "D"
$ws.Range("A6").Value = "E"
$ws.Range("A7").Value = "F"
$ws.Range("A8").Value = "G"
$ws.Range("A9").Value = "H"
$ws.Range("B2").Value = "Nenhuma Renda"
$ws.Range("B3").Value = "Até R$ 1.212,00"
$ws.Range("B4").Value = "De R$ 1.212,01 até R$ 1.818,00."
$ws.Range("B5").Value = "De R$ 1.818,01 até R$ 2.424,00."
$ws.Range("B6").Value = "De R$ 2.424,01 até R$ 3.030,00."
$ws.Range("B7").Value = "De R$ 3.030,01 até R$ 3.636,00."
$ws.Range("B8").Value = "De R$ 3.636,01 até R$ 4.848,00."
$ws.Range("B9").Value = "De R$ 4.848,01 até R$ 6.060,00."

$ws.Range("A10").Value = "I"
$ws.Range("B10").Value = "De R$ 6.060,01 até R$ 7.272,00."
$ws.Range("A11").Value = "J"
$ws.Range("B11").Value = "De R$ 7.272,01 até R$ 8.484,00."
$ws.Range("A12").Value = "K"
$ws.Range("B12").Value = "De R$ 8.484,01 até R$ 9.696,00."
$ws.Range("A13").Value = "L"
$ws.Range("B13").Value = "De R$ 9.696,01 até R$ 10.908,00."
$ws.Range("A14").Value = "M"
$ws.Range("B14").Value = "De R$ 10.908,01 até R$ 12.120,00."
$ws.Range("A15").Value = "N"
$ws.Range("B15").Value = "De R$ 12.120,01 até R$ 14.544,00."
$ws.Range("A16").Value = "O"
$ws.Range("B16").Value = "De R$ 14.544,01 até R$ 18.180,00."
$ws.Range("A17").Value = "P"
$ws.Range("B17").Value = "De R$ 18.180,01 até R$ 24.240,00."
$ws.Range("A18").Value = "Q"
$ws.Range("B18").Value = "Acima de R$ 24.240,00."

# Size the lookup columns to fit their contents, like the other DIME_* sheets.
$ws.Columns.Item(1).ColumnWidth = 17.21875
$ws.Columns.Item(2).ColumnWidth = 255.77734375

# Leave the new sheet as the active tab/sheet with the same selected cell
# recorded by the author (B31), and make sure it's the active window sheet.
$ws.Activate() | Out-Null
$ws.Range("B31").Select() | Out-Null
